$wb = $excel.ActiveWorkbook

# --- Sheet "Geral": refresh 46 team total scores (column B) ---
$geral = $wb.Worksheets.Item("Geral")
$geral.Range("B2").Value = 63.76
$geral.Range("B3").Value = 48.89
$geral.Range("B4").Value = 33.96
$geral.Range("B5").Value = 71.45999999999999
$geral.Range("B6").Value = 57.96
$geral.Range("B8").Value = 57.26
$geral.Range("B9").Value = 71.16
$geral.Range("B11").Value = 52.39
$geral.Range("B12").Value = 59.69
$geral.Range("B13").Value = 49.36
$geral.Range("B14").Value = 52.66
$geral.Range("B15").Value = 31.92
$geral.Range("B16").Value = 53.66
$geral.Range("B17").Value = 56.09
$geral.Range("B19").Value = 50.76
$geral.Range("B20").Value = 26.89
$geral.Range("B22").Value = 48.29
$geral.Range("B23").Value = 46.79
$geral.Range("B24").Value = 43.56
$geral.Range("B25").Value = 43.56
$geral.Range("B26").Value = 58.26
$geral.Range("B27").Value = 57.96
$geral.Range("B29").Value = 71.36
$geral.Range("B30").Value = 48.96
$geral.Range("B31").Value = 47.16
$geral.Range("B32").Value = 25.16
$geral.Range("B33").Value = 81.76000000000001
$geral.Range("B34").Value = 61.96
$geral.Range("B35").Value = 31.61
$geral.Range("B36").Value = 51.66
$geral.Range("B37").Value = 61.96
$geral.Range("B38").Value = 54.95
$geral.Range("B39").Value = 65.06
$geral.Range("B41").Value = 48.5
$geral.Range("B42").Value = 57.56
$geral.Range("B43").Value = 38.66
$geral.Range("B44").Value = 48.89
$geral.Range("B45").Value = 59.36
$geral.Range("B46").Value = 55.59
$geral.Range("B47").Value = 66.86
$geral.Range("B48").Value = 54.66
$geral.Range("B49").Value = 65.7
$geral.Range("B50").Value = 63.76
$geral.Range("B51").Value = 47.86
$geral.Range("B52").Value = 67.86
$geral.Range("B53").Value = 50.85

# --- Sheet "Mes - Janeiro": refresh ranking (re-sorted by new scores, desc) ---
$jan = $wb.Worksheets.Item("Mes - Janeiro")
$jan.Range("A2").Value = "Paulo Virgili FC"
$jan.Range("B2").Value = 81.76000000000001
$jan.Range("A3").Value = "C.A. Charru@"
$jan.Range("B3").Value = 71.70999999999999
$jan.Range("A4").Value = "Bandoleros FCS"
$jan.Range("B4").Value = 71.45999999999999
$jan.Range("A5").Value = "Mau Humor F.C."
$jan.Range("B5").Value = 71.36
$jan.Range("A6").Value = "dasdoresfc"
$jan.Range("B6").Value = 71.36
$jan.Range("A7").Value = "CARTOLEIRO DO VALLE PRO26.5"
$jan.Range("B7").Value = 71.16
$jan.Range("A8").Value = "Time do S.A.P.O"
$jan.Range("B8").Value = 67.86
$jan.Range("A9").Value = "Tatols Beants F.C"
$jan.Range("B9").Value = 66.86
$jan.Range("A10").Value = "teves_futsal20 f.c"
$jan.Range("B10").Value = 65.7
$jan.Range("A11").Value = "S.E.R. GRILLO"
$jan.Range("B11").Value = 65.06
$jan.Range("A12").Value = "A Lenda Super Vasco F.c"
$jan.Range("B12").Value = 63.76
$jan.Range("A13").Value = "Texas Club 2026"
$jan.Range("B13").Value = 63.76
$jan.Range("A14").Value = "PUXE FC"
$jan.Range("B14").Value = 61.96
$jan.Range("A15").Value = "Pity10"
$jan.Range("B15").Value = 61.96
$jan.Range("A16").Value = "Dom Camillo68"
$jan.Range("B16").Value = 59.69
$jan.Range("A17").Value = "Tabajara de Inhaua PB1"
$jan.Range("B17").Value = 59.36
$jan.Range("A18").Value = "MAFRA MARTINS FC"
$jan.Range("B18").Value = 58.51
$jan.Range("A19").Value = "lsauer fc"
$jan.Range("B19").Value = 58.26
$jan.Range("A20").Value = "BordonFC04"
$jan.Range("B20").Value = 57.96
$jan.Range("A21").Value = "Luis lemes inter"
$jan.Range("B21").Value = 57.96
$jan.Range("A22").Value = "Sport Clube PAIM"
$jan.Range("B22").Value = 57.56
$jan.Range("A23").Value = "cartola scheuer17"
$jan.Range("B23").Value = 57.26
$jan.Range("A24").Value = "Gremiomaniasm"
$jan.Range("B24").Value = 56.65
$jan.Range("A25").Value = "Fedato Futebol Clube"
$jan.Range("B25").Value = 56.09
$jan.Range("A26").Value = "TATITTA FC"
$jan.Range("B26").Value = 55.59
$jan.Range("A27").Value = "Rolo Compressor ZN"
$jan.Range("B27").Value = 54.95
$jan.Range("A28").Value = "TEAM LOPES 99"
$jan.Range("B28").Value = 54.66
$jan.Range("A29").Value = "FC Los Castilho"
$jan.Range("B29").Value = 53.66
$jan.Range("A30").Value = "FBC Colorado II"
$jan.Range("B30").Value = 52.66
$jan.Range("A31").Value = "DM Studio"
$jan.Range("B31").Value = 52.39
$jan.Range("A32").Value = "Profit Soccer"
$jan.Range("B32").Value = 51.66
$jan.Range("A33").Value = "VASCO MARTINS FC"
$jan.Range("B33").Value = 50.85
$jan.Range("A34").Value = "FÚRIA LEON"
$jan.Range("B34").Value = 50.76
$jan.Range("A35").Value = "FBC Colorado"
$jan.Range("B35").Value = 49.36
$jan.Range("A36").Value = "mercearia Estrela"
$jan.Range("B36").Value = 48.96
$jan.Range("A37").Value = "A Lenda Super Vascão f.c"
$jan.Range("B37").Value = 48.89
$jan.Range("A38").Value = "SUPER VASCÃO F.C"
$jan.Range("B38").Value = 48.89
$jan.Range("A39").Value = "SERGRILLO"
$jan.Range("B39").Value = 48.5
$jan.Range("A40").Value = "Grêmio imortal 37"
$jan.Range("B40").Value = 48.29
$jan.Range("A41").Value = "TIGRE LEON"
$jan.Range("B41").Value = 47.86
$jan.Range("A42").Value = "Máquina Laranjja"
$jan.Range("B42").Value = 47.16
$jan.Range("A43").Value = "JUV. KP"
$jan.Range("B43").Value = 46.79
$jan.Range("A44").Value = "FIGUEIRA DA ILHA"
$jan.Range("B44").Value = 44.06
$jan.Range("A45").Value = "LISI GREMISTA"
$jan.Range("B45").Value = 43.56
$jan.Range("A46").Value = "JV5 Tricolor Gaúcho"
$jan.Range("B46").Value = 43.56
$jan.Range("A47").Value = "Super Vasco f.c"
$jan.Range("B47").Value = 38.66
$jan.Range("A48").Value = "seralex"
$jan.Range("B48").Value = 34.36
$jan.Range("A49").Value = "AZURRA82"
$jan.Range("B49").Value = 33.96
$jan.Range("A50").Value = "FC castelo Branco 2"
$jan.Range("B50").Value = 31.92
$jan.Range("A51").Value = "pra sempre imortal fc"
$jan.Range("B51").Value = 31.61
$jan.Range("A52").Value = "Gig@ntte"
$jan.Range("B52").Value = 26.89
$jan.Range("A53").Value = "NaoVaiDescer!"
$jan.Range("B53").Value = 25.16
